$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected goldstandard link numbers
$ws.Range("I1").Value = 7610
$ws.Range("L1").Value = 1295
$ws.Range("O1").Value = 8240

# L1's original style ("s=1", explicit applyFont) reverts to the sheet's
# default/unset style (style index 0) when the corrected value is written -
# touching only NumberFormat (without forcing Font/Fill/Border/Alignment
# re-application) keeps the cell on the already-existing default xf record.
$ws.Range("L1").NumberFormat = "General"

# Move the selection to L1 (the corrected weight cell) as the new active cell
$ws.Range("L1").Select()

# Formulas in T4:V4 depend on C1/I1/L1/O1 and recompute automatically.
